# re-run RU 1001; without crop
# Updates a handful of recomputed mean values on Sheet 1.
# Columns: B = All, C = Europe, L = Russia (RU)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 55.3746767090015
$ws.Range("E2").Value = 56.7933697165467
$ws.Range("L2").Value = 49.0485989036895

$ws.Range("B3").Value = 40.563158268296
$ws.Range("L3").Value = 36.9618773883661

$ws.Range("B5").Value = 64.7210522905015

$ws.Range("B6").Value = 67.1354010141054
$ws.Range("C6").Value = 73.1501666372061
$ws.Range("L6").Value = 63.442019211072

$ws.Range("B7").Value = 68.4899159160604
$ws.Range("L7").Value = 60.1821932205212

$ws.Range("B8").Value = 61.8798692282585
$ws.Range("L8").Value = 54.0183622108344
